$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 'Java Integration Architect - W2 Contract - Ongoing'
$ws.Range("B14").Value = 'https://www.dice.com/job-detail/23ae0c27-91c5-4fc5-aa7c-90aaa609844f?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=3&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang'
$ws.Range("C14").Value = 'Hybrid in Oakland, California'
$ws.Range("D14").Value = 'Contract'
$ws.Range("E14").Value = '$75 - $85'
$ws.Range("F14").Value = 'Triune Infomatics Inc'
$ws.Range("G14").Value = 'Role: Java Integration Architect Location: Oakland, CA (2 days onsite mandatory) Duration: 3 months + extensions >> ongoing Overview: Lead a 3 months'' discovery, design, implementation support and remediation exercise so existing on-premises Java applications (Linux/JBoss) can interoperate with workloads migrating to cloud. Inspect code and runtime configuration to identify every change needed to repoint integrations and produce concise compatibility patterns to resolve common issues through testing and go-live. A major stream is re-platforming DB2 integrations to Microsoft SQL Server, covering both configuration and SQL/dialect changes. Project Insights: DB2 to SQL Server migration will involve confirmed complexities such as stored procedures; additional challenges like schema differences, SQL dialect variations, and data type mismatches are expected during delivery. Architect will collaborate with the existing DBA to analyze current models and resolve issues related to embedded SQL logic, pagination, identity/sequence handling, and stored procedure calls. Migration from on-prem DB2 to cloud SQL Server will require assessment and adaptation of TLS configurations, identity management practices, and network/proxy constraints — exact impacts to be identified as part of solution discovery. Key Responsibilities: Discovery & code/config review: Analyze Java apps on Linux/JBoss to locate all integration touchpoints (HTTP, messaging, databases, files, schedulers, security, network paths) with exact file/line, JNDI and environment references. Application Remediation Report (deliverable): Capture the current → target change set, change action, and validation steps for each Java application. Examples include: Connectivity : New endpoints/FQDNs/ports for migrated APIs/gateways/queues/databases; DNS/TTL; routing/peering/private link; proxy/NO_PROXY where required; DR/failover targets for those endpoints. Security : TLS specifics for new endpoints (server names/SANs, cert chains, trust/keystores, ciphers, hostname verification/SNI); identities/secrets used to call the migrated services (users/passwords, API keys/tokens), where stored, and rotation owner. Interfaces : HTTP APIs now fronted in cloud: base URLs, versions/paths, required headers/scopes, contract deltas and version pinning. Gateways referenced by Java (only if used): e.g., CTG/cloud API gateway host/port/alias/timeouts/pools/charset as called by the app. Messaging (JMS/MQ) only if that queueing endpoint is migrating: queue manager/channel/connection names, ports, TLS, credentials. Database : New driver and JDBC URL (host/port/db, encryption, failover/AG listener/multiSubnetFailover), credentials, and required SQL/dialect/schema changes (types, functions, pagination, isolation/locking/WITH UR equivalents, stored procedure call syntax, error/return code mapping, identifier/collation rules). Include the network path to SQL Server (private link/peering/proxy) and validation of connectivity/failover. Runtime & resilience : Connect/read timeouts, retries/back-off (idempotency where relevant), circuit breakers/health checks, transaction recovery impacts; charset/locale where gateway/database behavior changes. Application Compatibility Patterns (deliverable) : Produce short, reusable playbooks within Java scope for remediation and topology shifts (e.g., JDBC→SQL Server, TLS trust/hostname checks, HTTP client timeouts/retries, MQ/JMS client config, proxy/NO_PROXY, DNS failover, charset; and Java→CTG→CICS to Java→Cloud API Gateway→service; Java→VPN→DB2 to Java→SQL Server via private link), each with minimal code/config, verification and rollback. Stakeholder engagement: Work with client architects (application, integration, security, network) to shape options, agree patterns and support UAT/cutover. Governance & traceability: Maintain per-app change logs and risks; surface dependencies and non-functional impacts (latency, throughput, DR). Documentation: Design and document the implementation steps to guide the development team through execution. Update and revise the documentation to align with the final integration setup to train and hand over the maintainability of the integrations to the client support team. Skills & Experience (Essential): 7+ years as an application architect or senior Java engineer, able to read Java and pinpoint integration/configuration changes in Linux/JBoss environments. Demonstrable delivery of discovery/remediation for on-prem → cloud migrations, re-pointing integrations. Strong SQL with proven Microsoft SQL Server experience, including assessing/adapting SQL and schema when moving from DB2. Solid understanding in networking/DNS, routing/proxies and firewall/security-group models for app connectivity. Practical TLS/PKI knowledge for Java/JBoss (keystores/truststores, hostname verification, etc.). Clear technical writing and stakeholder management; comfortable operating independently while supporting client teams. Desirable: Familiarity with IBM MQ/JMS client configuration from Java. Exposure to feature-flag/remote-config approaches. Experience with common Java web stacks (Servlet/JAX-RS; Spring familiarity is a plus but not required). Experience with DR testing and performance tuning under changing latency profiles.'
$ws.Range("H14").Value = '35.79%%'
$ws.Range("I14").Value = 'Contract W2 | No Travel Required | Hybrid | $75 - $85/hr'
